$d = $word.ActiveDocument

# The "Goal:" paragraph previously read as one sentence about bridging
# neuroscience/AI. Replace it with the new project-goal sentence, keeping
# the existing (italic, Calibri) run formatting intact via MatchCase/Find-
# Replace rather than deleting+retyping the paragraph.
$old = "Goal: Bridging neuroscience and AI to decode mental representations and drive healthcare innovation."
$new = "Goal: Use neural data to improve AI models."

$find = $d.Content.Find
$ok = $find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)

Write-Output ("Replaced: " + $ok)
